$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.657.17"
$ws.Range("E2").Value = "  -1.22%  "
$ws.Range("D3").Value = "3.023.92"
$ws.Range("E3").Value = "  -1.44%  "
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").Value = "'582.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'148.63"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.03%  "
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("E8").Value = "  -2.78%  "
$ws.Range("D9").Value = "3.024.08"
$ws.Range("E9").Value = "  -1.23%  "
$ws.Range("E10").Value = "  -2.58%  "
$ws.Range("D11").Value = "'5.69"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.90%  "
$ws.Range("D12").Value = "'0.443"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.62%  "
$ws.Range("D13").Value = "'0.0000230"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.72%  "
$ws.Range("D14").Value = "'35.29"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.42%  "
$ws.Range("E15").Value = "  +2.02%  "
$ws.Range("D16").Value = "3.527.82"
$ws.Range("E16").Value = "  -1.54%  "
$ws.Range("D17").Value = "'7.06"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.68%  "
$ws.Range("D18").Value = "62.661.86"
$ws.Range("E18").Value = "  -1.44%  "
$ws.Range("D19").Value = "3.025.79"
$ws.Range("E19").Value = "  -1.69%  "
$ws.Range("D20").Value = "'467.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.33%  "
$ws.Range("D21").Value = "'14.03"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.39%  "
$ws.Range("D22").Value = "'0.691"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.09%  "
$ws.Range("D23").Value = "'7.41"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.85%  "
$ws.Range("D24").Value = "'2.36"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.16%  "
$ws.Range("D25").Value = "'80.89"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.21%  "
$ws.Range("D26").Value = "'12.38"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.90%  "
$ws.Range("D27").Value = "'10.37"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.25%  "
$ws.Range("D28").Value = "'0.998"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("D30").Value = "'7.24"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.59%  "
$ws.Range("E31").Value = "  -1.13%  "
$ws.Range("D32").Value = "'2.13"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.75%  "
$ws.Range("D33").Value = "'27.40"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.43%  "
$ws.Range("E34").Value = "  -4.09%  "
$ws.Range("D35").Value = "'1.03"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.89%  "
$ws.Range("D36").Value = "0.0₃0794"
$ws.Range("E36").Value = "  -2.73%  "
$ws.Range("E37").Value = "  -3.50%  "
$ws.Range("E38").Value = "  -2.15%  "
$ws.Range("D39").Value = "'50.22"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.56%  "
$ws.Range("E40").Value = "  -3.32%  "
$ws.Range("D41").Value = "'2.93"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -12.32%  "
$ws.Range("D42").Value = "'422.15"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.38%  "
$ws.Range("D43").Value = "'0.281"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.64%  "
$ws.Range("E44").Value = "  +0.62%  "
$ws.Range("D45").Value = "2.802.88"
$ws.Range("E45").Value = "  +0.22%  "
$ws.Range("D46").Value = "'0.0355"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.10%  "
$ws.Range("D47").Value = "'37.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.97%  "
$ws.Range("D48").Value = "'129.06"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.16%  "
$ws.Range("D50").Value = "'24.34"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.91%  "
$ws.Range("E51").Value = "  -0.82%  "
